# Apply updated market-price derived figures (H-N) across the per-job Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 74.57143000000001
$ws.Range("I6").Value = 13.4
$ws.Range("K6").Value = 40.2
$ws.Range("M6").Value = 71.8

$ws.Range("H11").Value = 1266.625
$ws.Range("I11").Value = 1266.625
$ws.Range("K11").Value = 1266.625
$ws.Range("M11").Value = -1126.625

$ws.Range("H17").Value = 6281.727
$ws.Range("J17").Value = 7472.1665
$ws.Range("L17").Value = 22416.4995
$ws.Range("N17").Value = -22752.4995

$ws.Range("H125").Value = 8774747
$ws.Range("I125").Value = 1350.6875
$ws.Range("K125").Value = 12156.1875
$ws.Range("M125").Value = -9696.1875

$ws.Range("H132").Value = 23257838
$ws.Range("I132").Value = 23257838
$ws.Range("K132").Value = 69773514
$ws.Range("M132").Value = -69770984

$ws.Range("H138").Value = 2881.024
$ws.Range("I138").Value = 1397.4584
$ws.Range("K138").Value = 4192.3752
$ws.Range("M138").Value = 947.6247999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 2362.5715
$ws.Range("I21").Value = 947
$ws.Range("J21").Value = 4250
$ws.Range("K21").Value = 947
$ws.Range("L21").Value = 4250
$ws.Range("M21").Value = -573
$ws.Range("N21").Value = -4998

$ws.Range("H36").Value = 5256
$ws.Range("I36").Value = 2703.1428
$ws.Range("J36").Value = 8830
$ws.Range("K36").Value = 2703.1428
$ws.Range("L36").Value = 8830
$ws.Range("M36").Value = -2357.1428
$ws.Range("N36").Value = -9522

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 15876014
$ws.Range("J99").Value = 3373.75
$ws.Range("L99").Value = 3373.75
$ws.Range("N99").Value = -6369.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 12000
$ws.Range("I4").Value = 12000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -11888
$ws.Range("N4").ClearContents()

$ws.Range("H31").Value = 27970.139
$ws.Range("J31").Value = 60592.312
$ws.Range("L31").Value = 60592.312
$ws.Range("N31").Value = -61182.312

$ws.Range("H34").Value = 27970.139
$ws.Range("J34").Value = 60592.312
$ws.Range("L34").Value = 60592.312
$ws.Range("N34").Value = -60996.312

$ws.Range("H58").Value = 8904.9375
$ws.Range("I58").Value = 14284.375
$ws.Range("J58").Value = 3525.5
$ws.Range("K58").Value = 14284.375
$ws.Range("L58").Value = 3525.5
$ws.Range("M58").Value = -14081.375
$ws.Range("N58").Value = -3931.5

$ws.Range("H136").Value = 8904.9375
$ws.Range("I136").Value = 14284.375
$ws.Range("J136").Value = 3525.5
$ws.Range("K136").Value = 42853.125
$ws.Range("L136").Value = 10576.5
$ws.Range("M136").Value = -40303.125
$ws.Range("N136").Value = -15676.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 189
$ws.Range("I6").Value = 189
$ws.Range("K6").Value = 567
$ws.Range("M6").Value = -454

$ws.Range("H54").Value = 4905
$ws.Range("J54").Value = 4905
$ws.Range("L54").Value = 14715
$ws.Range("N54").Value = -15833

$ws.Range("H56").Value = 16671913
$ws.Range("I56").Value = 16671913
$ws.Range("K56").Value = 16671913
$ws.Range("M56").Value = -16671383

$ws.Range("H68").Value = 1546.4445
$ws.Range("J68").Value = 2250.8
$ws.Range("L68").Value = 6752.400000000001
$ws.Range("N68").Value = -8374.400000000001

$ws.Range("H71").Value = 1546.4445
$ws.Range("J71").Value = 2250.8
$ws.Range("L71").Value = 20257.2
$ws.Range("N71").Value = -28369.2

$ws.Range("H97").Value = 472.42856
$ws.Range("I97").Value = 129
$ws.Range("J97").Value = 730
$ws.Range("K97").Value = 387
$ws.Range("L97").Value = 2190
$ws.Range("M97").Value = 109
$ws.Range("N97").Value = -3182

$ws.Range("H131").Value = 12629519
$ws.Range("I131").Value = 6412245
$ws.Range("J131").Value = 16670747
$ws.Range("K131").Value = 19236735
$ws.Range("L131").Value = 50012241
$ws.Range("M131").Value = -19231695
$ws.Range("N131").Value = -50022321

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18194700
$ws.Range("I70").Value = 25015750
$ws.Range("J70").Value = 5233.3335
$ws.Range("K70").Value = 25015750
$ws.Range("L70").Value = 5233.3335
$ws.Range("M70").Value = -25015480
$ws.Range("N70").Value = -5773.3335

$ws.Range("H73").Value = 18194700
$ws.Range("I73").Value = 25015750
$ws.Range("J73").Value = 5233.3335
$ws.Range("K73").Value = 25015750
$ws.Range("L73").Value = 5233.3335
$ws.Range("M73").Value = -25014814
$ws.Range("N73").Value = -7105.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7149.909
$ws.Range("I7").Value = 5219.4
$ws.Range("J7").Value = 8758.666999999999
$ws.Range("K7").Value = 5219.4
$ws.Range("L7").Value = 8758.666999999999
$ws.Range("M7").Value = -5107.4
$ws.Range("N7").Value = -8982.666999999999

$ws.Range("H46").Value = 3484.7715
$ws.Range("I46").Value = 1382.0555
$ws.Range("J46").Value = 5711.1763
$ws.Range("K46").Value = 1382.0555
$ws.Range("L46").Value = 5711.1763
$ws.Range("M46").Value = -1194.0555
$ws.Range("N46").Value = -6087.1763

$ws.Range("H55").Value = 1397.5416
$ws.Range("J55").Value = 1122.7693
$ws.Range("L55").Value = 1122.7693
$ws.Range("N55").Value = -1468.7693

$ws.Range("H68").Value = 499.5
$ws.Range("J68").Value = 499
$ws.Range("L68").Value = 499
$ws.Range("N68").Value = -1997

$ws.Range("H71").Value = 499.5
$ws.Range("J71").Value = 499
$ws.Range("L71").Value = 2495
$ws.Range("N71").Value = -9983

$ws.Range("H82").Value = 5053873
$ws.Range("I82").Value = 6947676
$ws.Range("K82").Value = 6947676
$ws.Range("M82").Value = -6947315

$ws.Range("H85").Value = 5053873
$ws.Range("I85").Value = 6947676
$ws.Range("K85").Value = 6947676
$ws.Range("M85").Value = -6946428

$ws.Range("H108").Value = 43626
$ws.Range("J108").Value = 43626
$ws.Range("L108").Value = 43626
$ws.Range("N108").Value = -51306

$ws.Range("H126").Value = 7149.909
$ws.Range("I126").Value = 5219.4
$ws.Range("J126").Value = 8758.666999999999
$ws.Range("K126").Value = 15658.2
$ws.Range("L126").Value = 26276.001
$ws.Range("M126").Value = -13188.2
$ws.Range("N126").Value = -31216.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 15502.25
$ws.Range("J31").Value = 18803.8
$ws.Range("L31").Value = 18803.8
$ws.Range("N31").Value = -19499.8

$ws.Range("H51").Value = 16321.7
$ws.Range("I51").Value = 9919.666999999999
$ws.Range("J51").Value = 25924.75
$ws.Range("K51").Value = 9919.666999999999
$ws.Range("L51").Value = 25924.75
$ws.Range("M51").Value = -9409.666999999999
$ws.Range("N51").Value = -26944.75

$ws.Range("H122").Value = 3358.111
$ws.Range("I122").Value = 3615.111
$ws.Range("J122").Value = 3101.111
$ws.Range("K122").Value = 10845.333
$ws.Range("L122").Value = 9303.332999999999
$ws.Range("M122").Value = -8395.332999999999
$ws.Range("N122").Value = -14203.333

$ws.Range("H126").Value = 3399.2727
$ws.Range("I126").Value = 3806.8572
$ws.Range("J126").Value = 2686
$ws.Range("K126").Value = 11420.5716
$ws.Range("L126").Value = 8058
$ws.Range("M126").Value = -8950.571599999999
$ws.Range("N126").Value = -12998
